$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the _GoBack bookmark that sits alone in its own paragraph
#    (right after "I plan to use the following data for this
#    analysis:"), turning that paragraph back into a plain empty one.
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2. Resize the table: explicit fixed width + new column widths.
#    NOTE: the Word object model expresses widths in points, while
#    OOXML stores dxa (twentieths of a point) -- divide by 20.
# ------------------------------------------------------------------
$tbl = $d.Tables(1)
$tbl.PreferredWidthType = 3   # wdPreferredWidthPoints (explicit width)
$tbl.PreferredWidth = 9175 / 20

$tbl.Columns(1).Width = 3685 / 20
$tbl.Columns(2).Width = 5490 / 20

# ------------------------------------------------------------------
# 3. Text edits inside the existing rows.
# ------------------------------------------------------------------
$d.Content.Find.Execute("Purpose", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Data to be extracted", 2) | Out-Null

$d.Content.Find.Execute("Venues in the city and comparison", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Venues in the city and comparison. Cluster of venues and their comparison", 2) | Out-Null

$d.Content.Find.Execute("Weather data", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Weather data for comparison", 2) | Out-Null

# ------------------------------------------------------------------
# 4. Append a new "Zillow" row to the table.
# ------------------------------------------------------------------
$newRow = $tbl.Rows.Add()
$tbl.Cell($newRow.Index, 1).Range.Text = "Zillow"
$cell2 = $tbl.Cell($newRow.Index, 2)
$cell2.Range.Text = "Range of rent in the city and average rent.X"

# Locate the boundary right after the final period (but this position
# coincides with the paragraph end, which the bookmark engine clamps
# to document start) -- so we bookmark the junk "X" boundary instead,
# then delete the "X", which leaves the bookmark correctly collapsed
# right after the period, matching the original placement style.
$cellStart = $cell2.Range.Start
$periodPos = $cellStart + "Range of rent in the city and average rent.".Length
$bmRange = $d.Range($periodPos, $periodPos)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
$d.Range($periodPos, $periodPos + 1).Delete()
